$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 19 (Testmail #17) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(19, 1).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item(19, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(19, 3).Value = "Testmail #17: Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item(19, 4).Value = "Planning / Afspraak"
$logs.Cells.Item(19, 5).Value = "Beste,`nBedank u voor uw e-mail. Ik zal een demo inplannen bij Van Dijk op vrijdag om 11:00 uur. `nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Cells.Item(19, 6).Value = "2025-07-31 21:59:14"
$logs.Cells.Item(19, 7).Value = "Ja"
$logs.Cells.Item(19, 8).Value = "Nee"
$logs.Cells.Item(19, 9).Value = "Ja"
$logs.Cells.Item(19, 10).Value = "Nee"

# Extend the conditional-formatting ranges on the Logs sheet from row 18 to row 19
$colLetters = @("D", "G", "H", "I", "J")
foreach ($col in $colLetters) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "18")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "19")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: append new row 7 (Planning / Afspraak, 1) ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(7, 1).Value = "Planning / Afspraak"
$dash.Cells.Item(7, 2).Value = 1

# Extend the bar chart's category/value series references from row 6 to row 7
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$7,'Dashboard'!`$B`$2:`$B`$7,1)"
